$d = $word.ActiveDocument

# 1. Update the "Sample Script Description" paragraph spacing: after 120 -> 240 (twips),
#    i.e. 6pt -> 12pt. This paragraph contains the text about missing lines.
$descText = "A few lines in the Sample Script are missing (Enter your code here). You need to complete the code as per the given instructions."
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $descText) {
        $p.Format.SpaceAfter = 12
        break
    }
}

# 2. Remove the duplicate "Sample Script:" subheader paragraph that immediately follows
#    the description paragraph (the first "Sample Script:" header stays).
$seenFirst = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq "Sample Script:") {
        if ($seenFirst) {
            $p.Range.Delete()
            break
        }
        $seenFirst = $true
    }
}
